# Applies the "Changes on Dec 29th" edit:
#  1. Splits a few runs and wraps the split point with w:proofErr
#     gramStart/gramEnd markers (mirrors Word's grammar-checker
#     placing a proofErr pair around the flagged span).
#  2. Splits the final "continue" run into "C" + "ontinue" (no proofErr).
#  3. Appends a handful of new paragraphs (incl. a Wingdings arrow demo)
#     after the last paragraph, before the section break.
#
# Because w:proofErr / w:sym are not reachable through the high-level
# Range/Find object model, every edit below is performed by replacing a
# paragraph's Range content with an explicit OOXML fragment via
# Range.InsertXML -- this lets us place proofErr markers and sym runs
# exactly where the canonical XML expects them, while leaving the
# paragraph's own pPr (style/numbering/etc.) untouched.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------
# 1. "Planning  (Intake)" -> "Planning  (" [gramStart/gramEnd] "Intake)"
# ---------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$xml1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Planning  (</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Intake)</w:t></w:r>
</w:p>
'@
$p1.Range.InsertXML($xml1)

# ---------------------------------------------------------------
# 2. "Feasibility Check  - RPA Management & RPA Solution " ->
#    "Feasibility " + [gramStart]"Check  -"[gramEnd] + " RPA Management & RPA Solution "
#    (rest of the paragraph - the "Archichet" spellcheck run and the
#    trailing "/Consultant - 100%" - is preserved unchanged)
# ---------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$xml3 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Feasibility </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Check  -</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> RPA Management &amp; RPA Solution </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Archichet</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>/Consultant &#8211; 100%</w:t></w:r>
</w:p>
'@
$p3.Range.InsertXML($xml3)

# ---------------------------------------------------------------
# 3. "For Each [sym] Defined set of array" ->
#    "For Each [sym] Defined set of " + [gramStart]"array"[gramEnd]
# ---------------------------------------------------------------
$p61 = $d.Paragraphs.Item(61)
$xml61 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">For Each </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:sym w:font="Wingdings" w:char="F0E0"/></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Defined set of </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>array</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
</w:p>
'@
$p61.Range.InsertXML($xml61)

# ---------------------------------------------------------------
# 4. "Continue : will stop the execution..." ->
#    [gramStart]"Continue :"[gramEnd] + " will stop the execution..."
# ---------------------------------------------------------------
$p71 = $d.Paragraphs.Item(71)
$xml71 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Continue :</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> will stop the execution of remaining steps and continue the loop</w:t></w:r>
</w:p>
'@
$p71.Range.InsertXML($xml71)

# ---------------------------------------------------------------
# 5. Final paragraph "continue" -> "C" + "ontinue" (plain run split,
#    no proofErr), then append the 7 new trailing paragraphs
#    (the divide/Divide Wingdings demo block) in the same call so
#    they land right after it and before the sectPr.
# ---------------------------------------------------------------
$p74 = $d.Paragraphs.Item($d.Paragraphs.Count)
$xml74 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>C</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>ontinue</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:pBdr><w:bottom w:val="single" w:sz="6" w:space="1" w:color="auto"/></w:pBdr>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-US"/></w:rPr></w:pPr>
</w:p>
<w:p>
  <w:pPr><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-US"/></w:rPr></w:pPr>
  <w:r>
    <w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-US"/></w:rPr>
    <w:t>&#8220;Divide&#8221; not equal to &#8220;divide&#8221;</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-US"/></w:rPr></w:pPr>
</w:p>
<w:p>
  <w:pPr><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-US"/></w:rPr></w:pPr>
  <w:r>
    <w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-US"/></w:rPr>
    <w:t xml:space="preserve">DIVIDE </w:t>
  </w:r>
  <w:r>
    <w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-US"/></w:rPr>
    <w:sym w:font="Wingdings" w:char="F0E0"/>
  </w:r>
  <w:r>
    <w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-US"/></w:rPr>
    <w:t xml:space="preserve"> DIVIDE</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-US"/></w:rPr></w:pPr>
  <w:r>
    <w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-US"/></w:rPr>
    <w:t xml:space="preserve">divide </w:t>
  </w:r>
  <w:r>
    <w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-US"/></w:rPr>
    <w:sym w:font="Wingdings" w:char="F0E0"/>
  </w:r>
  <w:r>
    <w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-US"/></w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-US"/></w:rPr>
    <w:t>divide</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
'@
$p74.Range.InsertXML($xml74)

Write-Output "done"
